$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted before the current row 211
# (everything from row 211 down shifts one row, 211->212, ..., 277->278).
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new record's data.
$ws.Cells.Item(211, 1).Value = 3
$ws.Cells.Item(211, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(211, 3).Value = "Coquimbo"
$ws.Cells.Item(211, 4).Value = 44588
$ws.Cells.Item(211, 5).Value = 5
$ws.Cells.Item(211, 6).Value = 100112009
$ws.Cells.Item(211, 7).Value = "Acelga"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 280
$ws.Cells.Item(211, 11).Value = 2300
$ws.Cells.Item(211, 12).Value = 2500
$ws.Cells.Item(211, 13).Value = 2386
$ws.Cells.Item(211, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(211, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(211, 16).Value = 398
$ws.Cells.Item(211, 17).Value = 6
$ws.Cells.Item(211, 18).Value = "Hortaliza"
